$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4306.467
$ws.Range("I51").Value = 3799.9
$ws.Range("J51").Value = 5319.6
$ws.Range("K51").Value = 3799.9
$ws.Range("L51").Value = 5319.6
$ws.Range("M51").Value = -3315.9
$ws.Range("N51").Value = -6287.6
# Row 69
$ws.Range("H69").Value = 34133.332
$ws.Range("I69").Value = 82333.336
$ws.Range("J69").Value = 22083.334
$ws.Range("K69").Value = 247000.008
$ws.Range("L69").Value = 66250.00199999999
$ws.Range("M69").Value = -246126.008
$ws.Range("N69").Value = -67998.00199999999
# Row 72
$ws.Range("H72").Value = 34133.332
$ws.Range("I72").Value = 82333.336
$ws.Range("J72").Value = 22083.334
$ws.Range("K72").Value = 741000.024
$ws.Range("L72").Value = 198750.006
$ws.Range("M72").Value = -736632.024
$ws.Range("N72").Value = -207486.006
# Row 76
$ws.Range("H76").Value = 3594.2856
$ws.Range("I76").Value = 3678.6667
$ws.Range("K76").Value = 3678.6667
$ws.Range("M76").Value = -3363.6667
# Row 79
$ws.Range("H79").Value = 3594.2856
$ws.Range("I79").Value = 3678.6667
$ws.Range("K79").Value = 3678.6667
$ws.Range("M79").Value = -2586.6667
# Row 80
$ws.Range("H80").Value = 2726.611
$ws.Range("I80").Value = 2709.5715
$ws.Range("K80").Value = 8128.7145
$ws.Range("M80").Value = -7130.7145
# Row 83
$ws.Range("H83").Value = 2726.611
$ws.Range("I83").Value = 2709.5715
$ws.Range("K83").Value = 24386.1435
$ws.Range("M83").Value = -19394.1435
# Row 99
$ws.Range("H99").Value = 251.4
$ws.Range("I99").Value = 277.33334
$ws.Range("K99").Value = 832.0000200000001
$ws.Range("M99").Value = 665.9999799999999
# Row 100
$ws.Range("H100").Value = 2415.1667
$ws.Range("I100").Value = 1623
$ws.Range("J100").Value = 3999.5
$ws.Range("K100").Value = 1623
$ws.Range("L100").Value = 3999.5
$ws.Range("M100").Value = -1082
$ws.Range("N100").Value = -5081.5
# Row 101
$ws.Range("H101").Value = 16667448
$ws.Range("I101").Value = 20000738
$ws.Range("K101").Value = 60002214
$ws.Range("M101").Value = -60000592
# Row 103
$ws.Range("H103").Value = 1048
$ws.Range("J103").Value = 1164
$ws.Range("L103").Value = 3492
$ws.Range("N103").Value = -4664
# Row 127
$ws.Range("H127").Value = 1839.3077
$ws.Range("I127").Value = 2187.8
$ws.Range("J127").Value = 1621.5
$ws.Range("K127").Value = 6563.400000000001
$ws.Range("L127").Value = 4864.5
$ws.Range("M127").Value = -1603.400000000001
$ws.Range("N127").Value = -14784.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 122.75
$ws.Range("I5").Value = 143.16667
$ws.Range("K5").Value = 143.16667
$ws.Range("M5").Value = -31.16667000000001
# Row 32
$ws.Range("H32").Value = 1751827.5
$ws.Range("I32").Value = 2060644.6
$ws.Range("J32").Value = 701849.2
$ws.Range("K32").Value = 2060644.6
$ws.Range("L32").Value = 701849.2
$ws.Range("M32").Value = -2060357.6
$ws.Range("N32").Value = -702423.2
# Row 97
$ws.Range("H97").Value = 1441.4286
$ws.Range("I97").Value = 1365.5
$ws.Range("J97").Value = 1897
$ws.Range("K97").Value = 1365.5
$ws.Range("L97").Value = 1897
$ws.Range("M97").Value = -869.5
$ws.Range("N97").Value = -2889

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 122.75
$ws.Range("I4").Value = 143.16667
$ws.Range("K4").Value = 143.16667
$ws.Range("M4").Value = -28.16667000000001
# Row 86
$ws.Range("H86").Value = 1450
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 89
$ws.Range("H89").Value = 1450
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 99
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
# Row 105
$ws.Range("H105").Value = 2484.75
$ws.Range("I105").Value = 2359.8572
$ws.Range("J105").Value = 2659.6
$ws.Range("K105").Value = 2359.8572
$ws.Range("L105").Value = 2659.6
$ws.Range("M105").Value = -612.8571999999999
$ws.Range("N105").Value = -6153.6
# Row 126
$ws.Range("H126").Value = 46999
$ws.Range("J126").Value = 46999
$ws.Range("L126").Value = 46999
$ws.Range("N126").Value = -56879
# Row 134
$ws.Range("H134").Value = 2680.3333
$ws.Range("I134").Value = 2680.3333
$ws.Range("K134").Value = 8040.999899999999
$ws.Range("M134").Value = -5505.999899999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2731.6667
$ws.Range("I16").Value = 3182.1667
$ws.Range("K16").Value = 3182.1667
$ws.Range("M16").Value = -2895.1667
# Row 58
$ws.Range("H58").Value = 2044.4
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797
# Row 113
$ws.Range("H113").Value = 2731.6667
$ws.Range("I113").Value = 3182.1667
$ws.Range("K113").Value = 3182.1667
$ws.Range("M113").Value = -1012.1667
# Row 136
$ws.Range("H136").Value = 2044.4
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450
# Row 141
$ws.Range("H141").Value = 108422.07
$ws.Range("J141").Value = 108422.07
$ws.Range("L141").Value = 108422.07
$ws.Range("N141").Value = -118782.07

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 69
$ws.Range("J33").Value = 69
$ws.Range("L33").Value = 414
$ws.Range("N33").Value = -980
# Row 75
$ws.Range("H75").Value = 5816.4287
$ws.Range("J75").Value = 7642.4
$ws.Range("L75").Value = 22927.2
$ws.Range("N75").Value = -24923.2
# Row 78
$ws.Range("H78").Value = 5816.4287
$ws.Range("J78").Value = 7642.4
$ws.Range("L78").Value = 68781.59999999999
$ws.Range("N78").Value = -78765.59999999999
# Row 92
$ws.Range("H92").Value = 633
$ws.Range("J92").Value = 599.5
$ws.Range("L92").Value = 1798.5
$ws.Range("N92").Value = -4294.5
# Row 122
$ws.Range("H122").Value = 100
$ws.Range("I122").Value = 100
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 900
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = 1550
$ws.Range("N122").Value = -5800
# Row 129
$ws.Range("H129").Value = 1114886.9
$ws.Range("J129").Value = 1670922
$ws.Range("L129").Value = 5012766
$ws.Range("N129").Value = -5022766
# Row 131
$ws.Range("H131").Value = 528779.7
$ws.Range("I131").Value = 1408
$ws.Range("K131").Value = 4224
$ws.Range("M131").Value = 816
# Row 139
$ws.Range("H139").Value = 571.1429000000001
$ws.Range("I139").Value = 400
$ws.Range("K139").Value = 1200
$ws.Range("M139").Value = 3940
# Row 140
$ws.Range("H140").Value = 11511.75
$ws.Range("I140").Value = 2998.75
$ws.Range("K140").Value = 8996.25
$ws.Range("M140").Value = -3816.25

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3083.3635
$ws.Range("I132").Value = 2891.9
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 8675.700000000001
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -6145.700000000001
$ws.Range("N132").Value = -20054

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1199.4
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
# Row 27
$ws.Range("H27").Value = 1199.4
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 93
$ws.Range("H93").Value = 3427.5
$ws.Range("I93").Value = 3078.75
$ws.Range("K93").Value = 3078.75
$ws.Range("M93").Value = -1830.75
# Row 100
$ws.Range("H100").Value = 3960.7144
$ws.Range("I100").Value = 3509
$ws.Range("J100").Value = 5090
$ws.Range("K100").Value = 3509
$ws.Range("L100").Value = 5090
$ws.Range("M100").Value = -2968
$ws.Range("N100").Value = -6172
# Row 132
$ws.Range("H132").Value = 6377.875
$ws.Range("J132").Value = 3916.3333
$ws.Range("L132").Value = 11748.9999
$ws.Range("N132").Value = -16808.9999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 6341.222
$ws.Range("I126").Value = 5224.067
$ws.Range("J126").Value = 7737.6665
$ws.Range("K126").Value = 15672.201
$ws.Range("L126").Value = 23212.9995
$ws.Range("M126").Value = -13202.201
$ws.Range("N126").Value = -28152.9995
